$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.972.57'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.015.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.79%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.39'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.51'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.699'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +12.22%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.747'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.78'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +7.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000324'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.70'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.657.94'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.86%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.050.22'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.02'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.48%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.18'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.011.89'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.57'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '97.50'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.50'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.18'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.30'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.12'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -8.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.71'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.85'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.67'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +18.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.61'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.43'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +7.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.42'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.131'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '680.40'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '48.16'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +18.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.17'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.445'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.151'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0821'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -9.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.38'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -8.29%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.33%  '

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0489'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.25'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +12.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.149'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.66'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.70%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.00'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.65%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.34'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.29%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000267'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.20%  '
